$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defect Report - Template")

# Update the "Incidence / Severity / Probability of reproduction" detail text (cell B13)
$ws.Range("B13").Value = "Critical severity - La funcion de crear usuario funciona incorrectamente."

# Update the zoom level of the active window/sheet view from 100% to 115%
$excel.ActiveWindow.Zoom = 115
